$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$ws.Range("M56").Select()
$p2 = $win.Panes.Item(2)
$p2.ScrollRow = 69
$p2.ScrollColumn = 1
Write-Host "Pane2 ScrollRow=$($p2.ScrollRow) ScrollColumn=$($p2.ScrollColumn)"
$p1 = $win.Panes.Item(1)
Write-Host "Pane1 ScrollRow=$($p1.ScrollRow) ScrollColumn=$($p1.ScrollColumn)"
